$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.098450064659119
$ws.Range("B1").Value = 1.582001209259033
$ws.Range("C1").Value = 4.614859580993652
$ws.Range("D1").Value = 0.4627927243709564
$ws.Range("E1").Value = 0.5145339369773865
